# Updates cryptos list (Price / Volume(1h) columns, plus a 3-row coin
# reshuffle at rows 41-43) to match the latest scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range('D2').Value = '58.394.68'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').Value = '2.356.61'
$ws.Range('E3').Value = '  +0.92%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '542.46'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.59'
$ws.Range('E6').Value = '  +2.76%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.34%  '
$ws.Range('E8').Value = '  +4.87%  '
$ws.Range('E9').Value = '  +0.75%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.60'
$ws.Range('E10').Value = '  +5.98%  '
$ws.Range('E11').Value = '  -0.55%  '
$ws.Range('E12').Value = '  +2.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '23.87'
$ws.Range('E13').Value = '  +1.94%  '
$ws.Range('D14').Value = '2.770.93'
$ws.Range('E14').Value = '  +0.84%  '
$ws.Range('D15').Value = '58.342.45'
$ws.Range('E15').Value = '  +2.22%  '
$ws.Range('E16').Value = '  +0.76%  '
$ws.Range('D17').Value = '2.331.87'
$ws.Range('E17').Value = '  -0.01%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '10.72'
$ws.Range('E18').Value = '  +2.91%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '333.18'
$ws.Range('E19').Value = '  -1.23%  '
$ws.Range('E20').Value = '  +2.99%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.77'
$ws.Range('E21').Value = '  -1.93%  '
$ws.Range('E22').Value = '  +0.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '62.97'
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('E24').Value = '  +1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '8.49'
$ws.Range('E25').Value = '  -2.90%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.998'
$ws.Range('E26').Value = '  +0.41%  '
$ws.Range('E27').Value = '  +6.22%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.76'
$ws.Range('E28').Value = '  +2.06%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.42'
$ws.Range('E29').Value = '  -0.06%  '
$ws.Range('D30').Value = '0.0₃0738'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.13'
$ws.Range('E31').Value = '  +0.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '18.43'
$ws.Range('E32').Value = '  -0.31%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.02'
$ws.Range('E33').Value = '  +12.27%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.26'
$ws.Range('E35').Value = '  +6.81%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.00'
$ws.Range('E36').Value = '  +0.57%  '
$ws.Range('E37').Value = '  -0.48%  '
$ws.Range('E38').Value = '  +4.71%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '39.20'
$ws.Range('E39').Value = '  +0.31%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '143.17'
$ws.Range('E40').Value = '  -3.66%  '
$ws.Range('B41').Value = 'Bittensor'
$ws.Range('C41').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '293.36'
$ws.Range('E41').Value = '  +5.15%  '
$ws.Range('B42').Value = 'PolygonEcosystemToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.378'
$ws.Range('E42').Value = '  +0.42%  '
$ws.Range('B43').Value = 'Filecoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.64'
$ws.Range('E43').Value = '  +1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0946'
$ws.Range('E44').Value = '  +2.15%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '19.23'
$ws.Range('E45').Value = '  +3.91%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0503'
$ws.Range('E46').Value = '  +0.50%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.564'
$ws.Range('E47').Value = '  +1.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0220'
$ws.Range('E48').Value = '  +1.73%  '
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.48'
$ws.Range('E50').Value = '  +1.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.08'
$ws.Range('E51').Value = '  +0.65%  '
